# Bump the example booking dates on the "Bookings" sheet from 2024 to 2026
# (same month/day, two years later):
#   B2: 2024-03-20 (45371) -> 2026-03-20 (46101)
#   B3: 2024-03-21 (45372) -> 2026-03-21 (46102)
# and move the saved cursor selection from C20 to B4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bookings")

$ws.Range("B2").Value = 46101
$ws.Range("B3").Value = 46102

$ws.Range("B4").Select()
